$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumValue($ref, $val) {
    $ws.Range($ref).Value2 = $val
}

function Set-NumStyled($ref, $donor, $val) {
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
    $ws.Range($ref).Value2 = $val
}

function Set-TextStyled($ref, $donor, $text) {
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# Header text updates (Volume number, date range)
$a8 = $ws.Range("A8").Value2
$idxVol = $a8.IndexOf("48")
$ws.Range("A8").Characters($idxVol + 1, 2).Text = "49"

$c9 = $ws.Range("C9").Value2
$idxD1 = $c9.IndexOf("11/27/2023")
$ws.Range("C9").Characters($idxD1 + 1, "11/27/2023".Length).Text = "12/4/2023"
$c9b = $ws.Range("C9").Value2
$idxD2 = $c9b.IndexOf("12/3/2023")
$ws.Range("C9").Characters($idxD2 + 1, "12/3/2023".Length).Text = "12/10/2023"

# Data table updates (rows 15-30)
Set-TextStyled "C15" "D15" "0"
Set-TextStyled "G15" "D15" "0"
Set-TextStyled "H15" "E15" "***.*"
Set-NumValue "L15" 40
Set-NumValue "N15" 0
Set-NumValue "C16" 3
Set-NumValue "D16" 7
Set-NumValue "E16" -57.142857142857
Set-NumValue "F16" 13
Set-NumValue "G16" 24
Set-NumValue "H16" -45.833333333333
Set-NumValue "I16" 214
Set-NumValue "J16" 202
Set-NumValue "K16" 5.940594059405
Set-NumValue "L16" 62.121212121212
Set-NumValue "M16" -16.078431372549
Set-NumValue "N16" -77.708333333333
Set-NumValue "C17" 3
Set-NumValue "D17" 4
Set-NumValue "E17" -25
Set-NumValue "F17" 17
Set-NumValue "G17" 17
Set-NumValue "H17" 0
Set-NumValue "I17" 266
Set-NumValue "J17" 285
Set-NumValue "K17" -6.666666666666
Set-NumValue "L17" 10.373443983402
Set-NumValue "M17" 27.272727272727
Set-NumValue "N17" -8.904109589041
Set-NumValue "C18" 3
Set-NumValue "D18" 7
Set-NumValue "E18" -57.142857142857
Set-NumValue "F18" 10
Set-NumValue "G18" 22
Set-NumValue "H18" -54.545454545454
Set-NumValue "I18" 185
Set-NumValue "J18" 258
Set-NumValue "K18" -28.294573643410
Set-NumValue "L18" -18.502202643171
Set-NumValue "M18" -56.572769953051
Set-NumValue "N18" -90.227152667723
Set-NumValue "C19" 16
Set-NumValue "D19" 7
Set-NumValue "E19" 128.571428571429
Set-NumValue "F19" 50
Set-NumValue "G19" 41
Set-NumValue "H19" 21.951219512195
Set-NumValue "I19" 657
Set-NumValue "J19" 610
Set-NumValue "K19" 7.704918032786
Set-NumValue "L19" 22.118959107806
Set-NumValue "M19" 60.635696821515
Set-NumValue "N19" 3.301886792452
Set-NumValue "C20" 7
Set-NumValue "D20" 7
Set-NumValue "E20" 0
Set-NumValue "F20" 26
Set-NumValue "G20" 31
Set-NumValue "H20" -16.129032258064
Set-NumValue "I20" 350
Set-NumValue "J20" 298
Set-NumValue "K20" 17.449664429530
Set-NumValue "L20" 73.267326732673
Set-NumValue "M20" 0
Set-NumValue "N20" -89.303178484107
Set-NumValue "C21" 32
Set-NumValue "D21" 32
Set-NumValue "E21" 0
Set-NumValue "F21" 118
Set-NumValue "G21" 135
Set-NumValue "H21" -12.592592592592
Set-NumValue "I21" 1697
Set-NumValue "J21" 1675
Set-NumValue "K21" 1.313432835820
Set-NumValue "L21" 25.055268975681
Set-NumValue "M21" 1.495215311004
Set-NumValue "N21" -76.071630005640
Set-NumStyled "C22" "C16" 1
Set-NumStyled "D22" "C16" 1
Set-NumStyled "E22" "E16" 0
Set-NumValue "G22" 3
Set-NumValue "H22" -66.666666666666
Set-NumValue "I22" 22
Set-NumValue "J22" 16
Set-NumValue "K22" 37.5
Set-NumValue "L22" 175
Set-NumValue "M22" 15.789473684210
Set-NumValue "C24" 26
Set-NumValue "D24" 33
Set-NumValue "E24" -21.212121212121
Set-NumValue "F24" 83
Set-NumValue "G24" 118
Set-NumValue "H24" -29.661016949152
Set-NumValue "I24" 1228
Set-NumValue "J24" 1401
Set-NumValue "K24" -12.348322626695
Set-NumValue "L24" -6.044376434583
Set-NumValue "M24" 22.922922922922
Set-NumValue "D25" 8
Set-NumValue "E25" 25
Set-NumValue "F25" 42
Set-NumValue "G25" 40
Set-NumValue "H25" 5
Set-NumValue "I25" 462
Set-NumValue "J25" 509
Set-NumValue "K25" -9.233791748526
Set-NumValue "L25" -1.702127659574
Set-NumValue "M25" -30.526315789473
Set-TextStyled "C26" "D15" "0"
Set-NumStyled "D26" "C16" 1
Set-NumStyled "E26" "E16" -100
Set-NumValue "G26" 1
Set-NumValue "H26" 200
Set-NumValue "J26" 32
Set-NumValue "K26" 3.125
Set-NumValue "L26" 57.142857142857
Set-TextStyled "C27" "D15" "0"
Set-NumValue "F27" 6
Set-NumValue "H27" 500
Set-NumValue "L27" 19.607843137254
Set-NumValue "N28" -80
Set-NumValue "N29" -81.481481481481
Set-NumStyled "C30" "C16" 1
Set-NumStyled "F30" "C16" 1
Set-NumValue "I30" 4
Set-NumValue "K30" -42.857142857142
Set-NumValue "L30" 0
